$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sql_statement")

# New data rows (lon, lat, gid, name) for rows 2..13 -- row 1 has the headers already.
$data = @(
    @(-5.11314713344928, 43.0475785813902, 1, "La Pandiela"),
    @(-5.11114634792837, 43.051480203257,  2, "Las Grandas"),
    @(-5.11548776300494, 43.0448571835119, 3, "Llano La Virgen"),
    @(-5.11642856729724, 43.0410265735243, 4, "El Castillo"),
    @(-5.11769350249134, 43.0487496631674, 5, "Mata Las Llanezas"),
    @(-5.11681136051063, 43.0469694049557, 6, "Las Bodegueras"),
    @(-5.12369887231031, 43.0469214270074, 7, "Prao San Pedro"),
    @(-5.11101103091742, 43.052865379369,  8, "Cota Real"),
    @(-5.11484621458358, 43.0494468327536, 9, "El Majadon"),
    @(-5.11249996871582, 43.045525626433,  10, "La Jagariz (Fuente)"),
    @(-5.11633289596547, 43.0434592693458, 11, "Tremapasaguas"),
    @(-5.11673010058391, 43.0428568786621, 12, "Prao Santiago")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
